$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('Cluster Name', 'Active cases'),
    @('3321 Rochester and Elmore District Health Service Yalunkan Aged Care Hostel Rochester', 13),
    @('3323 Villa Maria Catholic Homes St Bernadette''s Aged Care Sunshine North', 13),
    @('3600 Belvedere Aged Care Noble Park', 22),
    @('3601 Baptcare Westhaven community', 23),
    @('3653 Fronditha Thalpori St Albans Aged Care', 23),
    @('3939 Bupa Aged Care Eastwood', 11),
    @('4257 BlueCross The Gables Camberwell', 28),
    @('44087 Fitzroy Primary School Fitzroy', 10),
    @('44098 Stawell Primary School', 22),
    @('44121 Wallan Primary School Wallan', 14),
    @('44165 Greenvale Primary School', 20),
    @('44234 Lucknow Primary School Bairnsdale', 21),
    @('44667 Beaumaris Primary School Beaumaris', 20),
    @('44718 Parkdale Primary School Parkdale', 12),
    @('44811 Dandenong North Primary School Dandenong', 33),
    @('44812 Bairnsdale West Primary School', 10),
    @('44865 Parktone Primary School Parkdale', 28),
    @('44891 Cranbourne Park Primary School Cranbourne', 14),
    @('45248 Brookside P-9 College Caroline Springs', 19),
    @('45249 Creekside K-9 College Caroline Springs', 12),
    @('45267 Epping Views Primary School Epping', 13),
    @('45518 Ashwood High School Ashwood', 21),
    @('45569 Nhill College Nhill Outbreak', 37),
    @('45648 St Brendans Primary School Shepparton', 10),
    @('4574 Village Glen Aged Care Residences Mornington', 14),
    @('45784 Holy Rosary Primary School White Hills', 43),
    @('45836 St Joseph''s Primary School Sorrento', 10),
    @('46037 Nazareth Catholic Primary School Grovedale', 37),
    @('46050 Our Lady''s Catholic Primary School Craigieburn', 18),
    @('46125 Our Lady of the Southern Cross Primary School Manor Lakes', 20),
    @('46276 Hillcrest Christian College Clyde North', 15),
    @('46328 Ilim College Dallas Main Campus Dallas Oct', 19),
    @('51529 Sirius College Primary School Dallas', 10),
    @('52380 Al Iman College Melton South', 12),
    @('52473 John Henry Primary School Pakenham', 15),
    @('Alfred Health Caulfield Hospital Caulfield', 11),
    @('Camp Coolamatong Farm Camp Banksia Peninsula', 13),
    @('Escala NewQuay Construction Site Docklands Drive Docklands', 14),
    @('Gippsland and East Gippsland Aboriginal Co- Operative Bairnsdale', 12),
    @('Hamilton Country Music Festival Hamilton Golf Club Hamilton', 29),
    @('Hodges Real Estate Brighton', 11),
    @('Melton Willows Melton', 10),
    @('Oakleigh Grammar Melbourne Private School Oakleigh', 10),
    @('St Josephs Catholic Primary School Warragul', 13),
    @('Wagstaff Meat Processing Plant Cranbourne East', 34),
    @('Werribee Mercy Hospital Emergency Department', 11),
    @('Western Health Sunshine Hospital Emergency Department St Albans', 11)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$addr = $ws.UsedRange.Address()
Write-Host "Used range: $addr"